# edit.ps1
# Applies the "Updated cryptos list" data refresh to cryptos.xlsx.
#
# Strategy: for every changed cell we write the new text through a
# temporary formula ( ="<value>" ), then Copy + PasteSpecial(xlPasteValues)
# that single cell onto itself. This "bakes" the formula result down to a
# plain literal (string) cell - exactly like the source data - without
# leaving an <f> formula behind and, crucially, without Excel's normal
# Range.Value auto-type-detection turning numeric-looking text (e.g.
# "18.39", "1.24", "6.10") into real numbers (which would also silently
# re-stamp the cell's number format/style). PasteSpecial(-4163) is
# xlPasteValues.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-CellText "D2" '26.016.67'
Set-CellText "E2" '  -0.19%  '
Set-CellText "D3" '1.633.54'
Set-CellText "E4" '  +0.27%  '
Set-CellText "E5" '  -0.86%  '
Set-CellText "E6" '  -1.22%  '
Set-CellText "E7" '  +0.26%  '
Set-CellText "D8" '0.250'
Set-CellText "E8" '  -2.30%  '
Set-CellText "E9" '  -3.17%  '
Set-CellText "D10" '18.39'
Set-CellText "E10" '  -6.49%  '
Set-CellText "E11" '  -0.66%  '
Set-CellText "D12" '1.861.27'
Set-CellText "E12" '  -0.69%  '
Set-CellText "D13" '1.775.43'
Set-CellText "E13" '  +6.85%  '
Set-CellText "E14" '  -2.67%  '
Set-CellText "E15" '  -3.39%  '
Set-CellText "D16" '25.998.41'
Set-CellText "E16" '  -0.30%  '
Set-CellText "D17" '0.0₃0741'
Set-CellText "E17" '  -3.05%  '
Set-CellText "D18" '61.40'
Set-CellText "E18" '  -3.24%  '
Set-CellText "E19" '  +0.26%  '
Set-CellText "D20" '190.72'
Set-CellText "E20" '  -2.58%  '
Set-CellText "E21" '  -2.39%  '
Set-CellText "E22" '  -3.03%  '
Set-CellText "D23" '6.10'
Set-CellText "E23" '  -2.05%  '
Set-CellText "E24" '  +0.41%  '
Set-CellText "E25" '  -1.34%  '
Set-CellText "D26" '143.57'
Set-CellText "E26" '  -0.18%  '
Set-CellText "E27" '  +0.06%  '
Set-CellText "D28" '6.76'
Set-CellText "E28" '  -2.04%  '
Set-CellText "D29" '15.18'
Set-CellText "E29" '  -2.44%  '
Set-CellText "D30" '1.24'
Set-CellText "D31" '0.0482'
Set-CellText "E31" '  -3.37%  '
Set-CellText "E32" '  -4.53%  '
Set-CellText "D33" '3.11'
Set-CellText "E33" '  -5.43%  '
Set-CellText "B34" 'LidoDAOToken'
Set-CellText "C34" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText "D34" '1.49'
Set-CellText "E34" '  -3.19%  '
Set-CellText "B35" 'HuobiToken'
Set-CellText "C35" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText "D35" '2.40'
Set-CellText "E35" '  -2.31%  '
Set-CellText "D36" '1.132.29'
Set-CellText "E36" '  -0.12%  '
Set-CellText "D37" '0.861'
Set-CellText "E37" '  -5.15%  '
Set-CellText "E38" '  -0.97%  '
Set-CellText "E39" '  -4.67%  '
Set-CellText "E40" '  -1.66%  '
Set-CellText "D41" '98.42'
Set-CellText "E41" '  -0.98%  '
Set-CellText "E42" '  -2.90%  '
Set-CellText "B43" 'RocketPoolETH'
Set-CellText "C43" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-CellText "D43" '1.770.61'
Set-CellText "E43" '  -0.72%  '
Set-CellText "B44" 'FraxShare'
Set-CellText "C44" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText "D44" '5.22'
Set-CellText "E44" '  -5.13%  '
Set-CellText "E45" '  -1.48%  '
Set-CellText "D46" '54.77'
Set-CellText "E46" '  -3.49%  '
Set-CellText "D47" '0.0526'
Set-CellText "E47" '  -0.03%  '
Set-CellText "D48" '1.48'
Set-CellText "E48" '  +1.16%  '
Set-CellText "E49" '  -0.19%  '
Set-CellText "E50" '  +0.40%  '
Set-CellText "D51" '7.49'
Set-CellText "E51" '  -3.47%  '
